$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 58 to row 59 (keeps existing cellXfs, avoids creating new styles)
$ws.Range("A58:V58").Copy()
$ws.Range("A59:V59").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now set the values
$ws.Range("A59").Value = 58
$ws.Range("B59").Value = "malta"
$ws.Range("C59").Value = "premier-league"
$ws.Range("D59").Value = "2023-2024"
$ws.Range("E59").Value = 45262.67708333334
$ws.Range("F59").Value = "Balzan"
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = "Sirens"
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 1.61
$ws.Range("K59").Value = "01/12/2023 04:42"
$ws.Range("L59").Value = 1.95
$ws.Range("M59").Value = "02/12/2023 16:09"
$ws.Range("N59").Value = 3.6
$ws.Range("O59").Value = "01/12/2023 04:42"
$ws.Range("P59").Value = 3.1
$ws.Range("Q59").Value = "02/12/2023 16:09"
$ws.Range("R59").Value = 4.83
$ws.Range("S59").Value = "01/12/2023 04:42"
$ws.Range("T59").Value = 4.24
$ws.Range("U59").Value = "02/12/2023 16:09"
$ws.Range("V59").Value = "https://www.betexplorer.com/football/malta/premier-league/balzan-fc-sirens/byAwV0Xr/"

Write-Output "done"
